$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: C2 formula updated to average D2 with 1250 as well ---
$ws.Range("C2").Formula = "=AVERAGE(B2,1250)*(5/7)+AVERAGE(1250,D2)*(2/7)"

# --- Row 11: totals row reworked ---
$ws.Range("B11").Value = 0
$ws.Range("C11").Formula = "=(C3*C2)-SUM(C4:C10)"
$ws.Range("D11").Formula = "=(D3*D2)"

# --- New "Daily Per Item" mini table (rows 13-21), column A labels ---

# Header row 13
$ws.Range("A1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Daily Per Item"

$ws.Range("A1").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "Average"

# Row 14 - Milk (style copied from A4: bold, no-top box border)
$ws.Range("A4").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Milk"

# Row 15 - Cereal (style copied from A5: plain, left/right box border)
$ws.Range("A5").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Cereal"

# Rows 16-20 reuse the same left/right box border style as row 15
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Baby food"

$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Diapers"

$ws.Range("A15").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Bread"

$ws.Range("A15").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Peanut butter"

$ws.Range("A15").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Jelly/Jam"

# Row 21 - Other (totals row, full box border like A11)
$ws.Range("A11").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Other"

# --- Column B averages for the mini table ---

# Build the B14 style once (copy border/format from A4, then tweak to the
# centered integer look used by this little table), then propagate that
# finished format to the rest of the B14:B20 column so the same new style
# index gets reused instead of Excel minting one-off styles per cell.
$ws.Range("A4").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").NumberFormat = "0"
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("B14").VerticalAlignment = -4107
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").Formula = "=C4/6"

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Formula = "=C5/93"

$ws.Range("B14").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Formula = "=C6/162"

$ws.Range("B14").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Formula = "=C7/82"

$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Formula = "=C8/48"

$ws.Range("B14").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Formula = "=C9/20"

$ws.Range("B14").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Formula = "=C10/4"

# Row 21 total - full box border, same centered integer look
$ws.Range("A11").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").NumberFormat = "0"
$ws.Range("B21").HorizontalAlignment = -4108
$ws.Range("B21").VerticalAlignment = -4107
$ws.Range("B21").Font.Bold = $false
$ws.Range("B21").Formula = "=C11/(2075-415)"

$ws.Range("F17").Select()
